# Apply updated profit-calculation values across Sheets (scheduled runner refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 29974.25
$ws.Range("I21").Value = 39000
$ws.Range("J21").Value = 28971.389
$ws.Range("K21").Value = 39000
$ws.Range("L21").Value = 28971.389
$ws.Range("M21").Value = -38532
$ws.Range("N21").Value = -29907.389
$ws.Range("H23").Value = 29974.25
$ws.Range("I23").Value = 39000
$ws.Range("J23").Value = 28971.389
$ws.Range("K23").Value = 39000
$ws.Range("L23").Value = 28971.389
$ws.Range("M23").Value = -38766
$ws.Range("N23").Value = -29439.389

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 28504.5
$ws.Range("I30").Value = 28504.5
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 28504.5
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -28354.5
$ws.Range("N30").ClearContents()
$ws.Range("H61").Value = 3169.7585
$ws.Range("I61").Value = 1956.4445
$ws.Range("K61").Value = 1956.4445
$ws.Range("M61").Value = -1744.4445
$ws.Range("H74").Value = 5447.6333
$ws.Range("I74").Value = 1580
$ws.Range("K74").Value = 1580
$ws.Range("M74").Value = -706
$ws.Range("H77").Value = 5447.6333
$ws.Range("I77").Value = 1580
$ws.Range("K77").Value = 7900
$ws.Range("M77").Value = -3532
$ws.Range("H122").Value = 2615
$ws.Range("I122").Value = 2012
$ws.Range("J122").Value = 2916.5
$ws.Range("K122").Value = 6036
$ws.Range("L122").Value = 8749.5
$ws.Range("M122").Value = -3586
$ws.Range("N122").Value = -13649.5
$ws.Range("H132").Value = 4365.909
$ws.Range("I132").Value = 4592.25
$ws.Range("J132").Value = 4236.5713
$ws.Range("K132").Value = 13776.75
$ws.Range("L132").Value = 12709.7139
$ws.Range("M132").Value = -11246.75
$ws.Range("N132").Value = -17769.7139
$ws.Range("H136").Value = 3169.7585
$ws.Range("I136").Value = 1956.4445
$ws.Range("K136").Value = 5869.333500000001
$ws.Range("M136").Value = -3319.333500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 508
$ws.Range("I29").Value = 508
$ws.Range("K29").Value = 508
$ws.Range("M29").Value = -219
$ws.Range("H94").Value = 2104.75
$ws.Range("I94").Value = 1904.5
$ws.Range("J94").Value = 2305
$ws.Range("K94").Value = 1904.5
$ws.Range("L94").Value = 2305
$ws.Range("M94").Value = -1453.5
$ws.Range("N94").Value = -3207
$ws.Range("H134").Value = 3029.6316
$ws.Range("I134").Value = 1964.3043
$ws.Range("K134").Value = 5892.9129
$ws.Range("M134").Value = -3357.9129

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2993.5186
$ws.Range("I58").Value = 2139.6365
$ws.Range("J58").Value = 3580.5625
$ws.Range("K58").Value = 2139.6365
$ws.Range("L58").Value = 3580.5625
$ws.Range("M58").Value = -1936.6365
$ws.Range("N58").Value = -3986.5625
$ws.Range("H107").Value = 430.1875
$ws.Range("I107").Value = 194.54546
$ws.Range("J107").Value = 948.6
$ws.Range("K107").Value = 194.54546
$ws.Range("L107").Value = 948.6
$ws.Range("M107").Value = 1725.45454
$ws.Range("N107").Value = -4788.6
$ws.Range("H132").Value = 3890.2
$ws.Range("I132").Value = 3635.1428
$ws.Range("J132").Value = 4027.5386
$ws.Range("K132").Value = 10905.4284
$ws.Range("L132").Value = 12082.6158
$ws.Range("M132").Value = -8375.428400000001
$ws.Range("N132").Value = -17142.6158
$ws.Range("H134").Value = 3548.158
$ws.Range("I134").Value = 1645.25
$ws.Range("J134").Value = 4932.091
$ws.Range("K134").Value = 4935.75
$ws.Range("L134").Value = 14796.273
$ws.Range("M134").Value = -2400.75
$ws.Range("N134").Value = -19866.273
$ws.Range("H136").Value = 2993.5186
$ws.Range("I136").Value = 2139.6365
$ws.Range("J136").Value = 3580.5625
$ws.Range("K136").Value = 6418.9095
$ws.Range("L136").Value = 10741.6875
$ws.Range("M136").Value = -3868.9095
$ws.Range("N136").Value = -15841.6875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2330.2856
$ws.Range("I70").Value = 1462.4
$ws.Range("J70").Value = 4500
$ws.Range("K70").Value = 4387.200000000001
$ws.Range("L70").Value = 13500
$ws.Range("M70").Value = -4072.200000000001
$ws.Range("N70").Value = -14130
$ws.Range("H73").Value = 2330.2856
$ws.Range("I73").Value = 1462.4
$ws.Range("J73").Value = 4500
$ws.Range("K73").Value = 4387.200000000001
$ws.Range("L73").Value = 13500
$ws.Range("M73").Value = -3295.200000000001
$ws.Range("N73").Value = -15684
$ws.Range("H131").Value = 1505.661
$ws.Range("I131").Value = 820
$ws.Range("J131").Value = 1542.3928
$ws.Range("K131").Value = 2460
$ws.Range("L131").Value = 4627.178400000001
$ws.Range("M131").Value = 2580
$ws.Range("N131").Value = -14707.1784

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1390941.6
$ws.Range("I122").Value = 2223585
$ws.Range("J122").Value = 3202.6667
$ws.Range("K122").Value = 6670755
$ws.Range("L122").Value = 9608.000100000001
$ws.Range("M122").Value = -6668305
$ws.Range("N122").Value = -14508.0001
$ws.Range("H132").Value = 3053.0667
$ws.Range("I132").Value = 2558.2593
$ws.Range("J132").Value = 3795.2778
$ws.Range("K132").Value = 7674.777900000001
$ws.Range("L132").Value = 11385.8334
$ws.Range("M132").Value = -5144.777900000001
$ws.Range("N132").Value = -16445.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 377.72726
$ws.Range("I55").Value = 296.125
$ws.Range("J55").Value = 595.3333
$ws.Range("K55").Value = 296.125
$ws.Range("L55").Value = 595.3333
$ws.Range("M55").Value = -123.125
$ws.Range("N55").Value = -941.3333
$ws.Range("H136").Value = 4865.8276
$ws.Range("I136").Value = 2639.1333
$ws.Range("J136").Value = 7251.5713
$ws.Range("K136").Value = 7917.3999
$ws.Range("L136").Value = 21754.7139
$ws.Range("M136").Value = -5367.3999
$ws.Range("N136").Value = -26854.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 125001450
$ws.Range("J96").Value = 1407.3334
$ws.Range("L96").Value = 1407.3334
$ws.Range("N96").Value = -4153.3334
$ws.Range("H132").Value = 19233554
$ws.Range("I132").Value = 25002252
$ws.Range("K132").Value = 75006756
$ws.Range("M132").Value = -75004226
$ws.Range("H136").Value = 13375428
$ws.Range("I136").Value = 37149164
$ws.Range("K136").Value = 111447492
$ws.Range("M136").Value = -111444942

